# Generate Report for Handback
# Update status/handback info on the Overview, zh-cn and de-de sheets
# to reflect that the handback is complete and in sync with en-US.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: status columns for zh-cn (E2) and de-de (F2)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# zh-cn sheet: status, latest handback datetime, error detail (now cleared)
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-09-03 08:53:15"
$p2ZhCnStyle = $wsZhCn.Range("P2").Style
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P2").Style = $p2ZhCnStyle

# de-de sheet: status, latest handback datetime, error detail (now cleared)
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-09-03 08:53:22"
$p2DeDeStyle = $wsDeDe.Range("P2").Style
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P2").Style = $p2DeDeStyle

# Re-autofit the columns whose text changed: the Status columns grew wider
# (longer text) and the Error Detail columns shrank (now empty).
$wsOverview.Columns.Item(5).ColumnWidth = 29.1   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 29.1   # F: de-de status

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1       # C: Status
$wsZhCn.Columns.Item(16).ColumnWidth = 12.8      # P: Error Detail

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1       # C: Status
$wsDeDe.Columns.Item(16).ColumnWidth = 12.8      # P: Error Detail
